$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: name and card number changes
$ws.Range("C2").Value = "Hartmut"

# Card number must remain text (not be coerced to a number, and keep the
# original cell style). Assigning the digit string directly via .Value
# would auto-convert it to a number, so we build it as a text formula
# result first, then paste just the value back over itself - this keeps
# the original style index and produces a plain text cell.
$ws.Range("B3").Formula = "=""2570314725427075"""
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# Row 5: statement start balance date
$ws.Range("D5").Value = "KONTOSTAND AM 23.03.2024"

# Row 6
$ws.Range("B6").Value = "25.03."
$ws.Range("C6").Value = "26.03."
$ws.Range("D6").Value = "PAYPAL YULXRI"
$ws.Range("E6").Value = "72,62-"

# Row 7
$ws.Range("B7").Value = "28.03."
$ws.Range("C7").Value = "29.03."
$ws.Range("D7").Value = "PAYPAL AJKVRS"
$ws.Range("E7").Value = "24,82-"

# Row 8
$ws.Range("B8").Value = "31.03."
$ws.Range("C8").Value = "01.04."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 58234010"
$ws.Range("E8").Value = "86,26-"

# Row 9
$ws.Range("B9").Value = "04.04."
$ws.Range("C9").Value = "05.04."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 95129290"
$ws.Range("E9").Value = "40,07-"

# Row 10: newly filled in (was empty), copy E column formatting from the row above
# (E10 previously used the "blank row" style; data rows use the right-aligned style)
$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B10").Value = "05.04."
$ws.Range("C10").Value = "06.04."
$ws.Range("D10").Value = "PAYPAL UZMMSE"
$ws.Range("E10").Value = "42,38-"

# Row 12: final statement balance
$ws.Range("D12").Value = "KONTOSTAND AM 10.04.2024"
$ws.Range("E12").Value = "266,15-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 19.04.2024"
